$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-bucket the Department column (B) for the rows whose Order Status (A)
# relates to the auto packer, CS confirmation, or holds at the various
# stores, into three new, more specific department values.

$autoPackerRows = 10, 140
foreach ($r in $autoPackerRows) {
    $ws.Cells.Item($r, 2).Value = "Auto Packer"
}

# Row 4 keeps its default (unstyled) formatting when retyped.
$ws.Cells.Item(4, 2).Value = "CONFIRMATION"

# These rows were previously styled (left-aligned, style index 2) and lose
# that explicit formatting when retyped with the new department value.
$confirmationRowsRestyled = 41, 42, 43, 44, 45, 46, 171, 172, 173, 174, 175, 176
foreach ($r in $confirmationRowsRestyled) {
    $c = $ws.Cells.Item($r, 2)
    $c.Value = "CONFIRMATION"
    $c.Style = "Normal"
}

$holdRows = 6, 72, 73, 74, 202, 203, 204
foreach ($r in $holdRows) {
    $ws.Cells.Item($r, 2).Value = "HOLD"
}

$ws.Range("A6").Select()

# The sheet grew from 268 to 289 data rows since the filter was last set;
# refresh the hidden _FilterDatabase defined name to match.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Sheet1!`$A`$1:`$C`$289"
